$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.479.19"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.617.88"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.33"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.84"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.847.56"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "1.623.50"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.92"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "27.462.45"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.59"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.13"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +5.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.55"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.85"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "1.472.57"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.954"
$ws.Range("E37").Value = "  +7.51%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.558"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0167"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.861"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.04"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.982"
$ws.Range("E43").Value = "  -4.91%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.28"
$ws.Range("E45").Value = "  -7.75%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.758.12"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.67"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.71"
$ws.Range("E51").Value = "  -1.13%  "
